# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" everywhere it appears
# - Narrower "Status" column(s) to fit the new (shorter) text

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: zh-cn / de-de status columns (E, F) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewRange = $wsOverview.Range("E2:F4")
foreach ($cell in $overviewRange.Cells) {
    if ($cell.Text -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- zh-cn sheet: Status column (C) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhRange = $wsZh.Range("C2:C4")
foreach ($cell in $zhRange.Cells) {
    if ($cell.Text -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- de-de sheet: Status column (C) ---
$wsDe = $wb.Worksheets.Item("de-de")
$deRange = $wsDe.Range("C2:C4")
foreach ($cell in $deRange.Cells) {
    if ($cell.Text -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

# --- Shrink the status columns now that the text is shorter ---
$wsOverview.Columns(5).ColumnWidth = 12.43
$wsOverview.Columns(6).ColumnWidth = 12.43
$wsZh.Columns(3).ColumnWidth = 12.43
$wsDe.Columns(3).ColumnWidth = 12.43
